# Session 3 update to the Minaria Cheatsheet workbook.
#  1. Two new header cells are added to the right of the existing table
#     (L2 = "Events", M2 = "Allignments"), matching the formatting of the
#     plain (non-header-shaded) cells already used elsewhere on the sheet.
#  2. Two Genasi names in column D are renamed for the new session:
#       "Marble"  -> "Jade"
#       "Whistle" -> "Tumult"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New columns for this session -------------------------------------
$ws.Range("L2").Value = "Events"
$ws.Range("M2").Value = "Allignments"

# Match the look of existing plain text cells (no header shading) instead of
# leaving the new cells with the default style - copy formatting only from a
# plain, unshaded data cell already on the sheet.
$ws.Range("B3").Copy()
$ws.Range("L2:M2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 2. Text corrections from the session -------------------------------
$ws.Range("D16").Value = "Jade"
$ws.Range("D22").Value = "Tumult"
